$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.033094729696814
$ws.Cells.Item(2, 4).Value = 1.035023211189406
$ws.Cells.Item(2, 5).Value = 0.9926147277508489
$ws.Cells.Item(2, 6).Value = 1.041687953765123
$ws.Cells.Item(2, 9).Value = 1.034428272690362
$ws.Cells.Item(2, 10).Value = 1.038221341237811
$ws.Cells.Item(2, 11).Value = 1.037821005770137
$ws.Cells.Item(2, 12).Value = 0.9955398523336033
$ws.Cells.Item(2, 13).Value = 1.04446676183946
$ws.Cells.Item(2, 14).Value = 1.039695734628047

# Row 3
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.034770753831877
$ws.Cells.Item(3, 4).Value = 1.036266888631019
$ws.Cells.Item(3, 5).Value = 0.9936372048519304
$ws.Cells.Item(3, 6).Value = 1.043452512421795
$ws.Cells.Item(3, 9).Value = 1.0348928722299
$ws.Cells.Item(3, 10).Value = 1.039536252295286
$ws.Cells.Item(3, 11).Value = 1.038872836829243
$ws.Cells.Item(3, 12).Value = 0.9963617723202692
$ws.Cells.Item(3, 13).Value = 1.046039498175515
$ws.Cells.Item(3, 14).Value = 1.04101251301002

# Row 4
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.035851639967672
$ws.Cells.Item(4, 4).Value = 1.037068344143259
$ws.Cells.Item(4, 5).Value = 0.9942998659930995
$ws.Cells.Item(4, 6).Value = 1.044591067257671
$ws.Cells.Item(4, 9).Value = 1.035190292103258
$ws.Cells.Item(4, 10).Value = 1.040383208032958
$ws.Cells.Item(4, 11).Value = 1.039549610107224
$ws.Cells.Item(4, 12).Value = 0.9968940712668345
$ws.Cells.Item(4, 13).Value = 1.047053498422462
$ws.Cells.Item(4, 14).Value = 1.041860671521987

# Row 5
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.036305195437296
$ws.Cells.Item(5, 4).Value = 1.037404500177347
$ws.Cells.Item(5, 5).Value = 0.9945786998346017
$ws.Cells.Item(5, 6).Value = 1.045068958123512
$ws.Cells.Item(5, 9).Value = 1.035314564967746
$ws.Cells.Item(5, 10).Value = 1.040738352260194
$ws.Cells.Item(5, 11).Value = 1.039833218674058
$ws.Cells.Item(5, 12).Value = 0.997117960005301
$ws.Cells.Item(5, 13).Value = 1.047478921644934
$ws.Cells.Item(5, 14).Value = 1.04221632009476

# Row 6
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.036381300069715
$ws.Cells.Item(6, 4).Value = 1.037460897076179
$ws.Cells.Item(6, 5).Value = 0.9946255319796338
$ws.Cells.Item(6, 6).Value = 1.045149154148117
$ws.Cells.Item(6, 9).Value = 1.035335386371866
$ws.Cells.Item(6, 10).Value = 1.040797929196257
$ws.Cells.Item(6, 11).Value = 1.039880784918519
$ws.Cells.Item(6, 12).Value = 0.9971555583673453
$ws.Cells.Item(6, 13).Value = 1.047550301932708
$ws.Cells.Item(6, 14).Value = 1.042275981636906

# Row 7
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.035857703712411
$ws.Cells.Item(7, 4).Value = 1.037072838910112
$ws.Cells.Item(7, 5).Value = 0.9943035907978915
$ws.Cells.Item(7, 6).Value = 1.044597455813352
$ws.Cells.Item(7, 9).Value = 1.035191955631931
$ws.Cells.Item(7, 10).Value = 1.040387957069077
$ws.Cells.Item(7, 11).Value = 1.039553403248089
$ws.Cells.Item(7, 12).Value = 0.9968970624459041
$ws.Cells.Item(7, 13).Value = 1.047059186316851
$ws.Cells.Item(7, 14).Value = 1.041865427302282

# Row 8
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.033661907805316
$ws.Cells.Item(8, 4).Value = 1.03544420506921
$ws.Cells.Item(8, 5).Value = 0.9929600610674294
$ws.Cells.Item(8, 6).Value = 1.042284973819493
$ws.Cells.Item(8, 9).Value = 1.034585953559973
$ws.Cells.Item(8, 10).Value = 1.038666532652597
$ws.Cells.Item(8, 11).Value = 1.038177276688953
$ws.Cells.Item(8, 12).Value = 0.9958175282591053
$ws.Cells.Item(8, 13).Value = 1.044999043454808
$ws.Cells.Item(8, 14).Value = 1.04014155826568

# Row 9
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029764235838083
$ws.Cells.Item(9, 4).Value = 1.03254867113996
$ws.Cells.Item(9, 5).Value = 0.9906006454969559
$ws.Cells.Item(9, 6).Value = 1.038184591162249
$ws.Cells.Item(9, 9).Value = 1.033493296518491
$ws.Cells.Item(9, 10).Value = 1.035602870036975
$ws.Cells.Item(9, 11).Value = 1.035722548853969
$ws.Cells.Item(9, 12).Value = 0.9939188001724441
$ws.Cells.Item(9, 13).Value = 1.041340058952694
$ws.Cells.Item(9, 14).Value = 1.037073544897737

# Row 10
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.027145676352406
$ws.Cells.Item(10, 4).Value = 1.030600351262787
$ws.Cells.Item(10, 5).Value = 0.989033133672735
$ws.Cells.Item(10, 6).Value = 1.035432826841936
$ws.Cells.Item(10, 9).Value = 1.032747840857804
$ws.Cells.Item(10, 10).Value = 1.033539245440912
$ws.Cells.Item(10, 11).Value = 1.03406535373186
$ws.Cells.Item(10, 12).Value = 0.9926553831429383
$ws.Cells.Item(10, 13).Value = 1.03888045308385
$ws.Cells.Item(10, 14).Value = 1.035006989718047

# Row 11
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.026006812893266
$ws.Cells.Item(11, 4).Value = 1.029752288524085
$ws.Cells.Item(11, 5).Value = 0.988355674866747
$ws.Cells.Item(11, 6).Value = 1.034236739362818
$ws.Cells.Item(11, 9).Value = 1.032420937170675
$ws.Cells.Item(11, 10).Value = 1.032640465857694
$ws.Cells.Item(11, 11).Value = 1.033342708972075
$ws.Cells.Item(11, 12).Value = 0.9921088820399291
$ws.Cells.Item(11, 13).Value = 1.037810391950828
$ws.Cells.Item(11, 14).Value = 1.034106933764728

# Row 12
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.025583015825673
$ws.Cells.Item(12, 4).Value = 1.029436601360021
$ws.Cells.Item(12, 5).Value = 0.9881042295826724
$ws.Cells.Item(12, 6).Value = 1.033791754697762
$ws.Cells.Item(12, 9).Value = 1.032298885416121
$ws.Cells.Item(12, 10).Value = 1.032305819471107
$ws.Cells.Item(12, 11).Value = 1.033073512114785
$ws.Cells.Item(12, 12).Value = 0.9919059725120875
$ws.Cells.Item(12, 13).Value = 1.037412148506363
$ws.Cells.Item(12, 14).Value = 1.033771812141888

# Row 13
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.025673957011573
$ws.Cells.Item(13, 4).Value = 1.029504348272645
$ws.Cells.Item(13, 5).Value = 0.9881581567098651
$ws.Cells.Item(13, 6).Value = 1.033887237625814
$ws.Cells.Item(13, 9).Value = 1.032325094330866
$ws.Cells.Item(13, 10).Value = 1.032377638699612
$ws.Cells.Item(13, 11).Value = 1.033131291008358
$ws.Cells.Item(13, 12).Value = 0.9919494934313052
$ws.Cells.Item(13, 13).Value = 1.037497608452266
$ws.Cells.Item(13, 14).Value = 1.033843733361936

# Row 14
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.025971797569
$ws.Cells.Item(14, 4).Value = 1.029726207654176
$ws.Cells.Item(14, 5).Value = 0.9883348863814464
$ws.Cells.Item(14, 6).Value = 1.034199971280934
$ws.Cells.Item(14, 9).Value = 1.032410861127038
$ws.Cells.Item(14, 10).Value = 1.032612820287994
$ws.Cells.Item(14, 11).Value = 1.033320472949335
$ws.Cells.Item(14, 12).Value = 0.9920921077337197
$ws.Cells.Item(14, 13).Value = 1.037777488943408
$ws.Cells.Item(14, 14).Value = 1.034079248935149

# Row 15
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.026155204010217
$ws.Cells.Item(15, 4).Value = 1.029862812166019
$ws.Cells.Item(15, 5).Value = 0.9884438009545853
$ws.Cells.Item(15, 6).Value = 1.034392562845135
$ws.Cells.Item(15, 9).Value = 1.032463621837193
$ws.Cells.Item(15, 10).Value = 1.032757616988177
$ws.Cells.Item(15, 11).Value = 1.033436931178509
$ws.Cells.Item(15, 12).Value = 0.9921799884222134
$ws.Cells.Item(15, 13).Value = 1.037949829284238
$ws.Cells.Item(15, 14).Value = 1.034224251263255

# Row 16
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.027221151732975
$ws.Cells.Item(16, 4).Value = 1.0306565399454
$ws.Cells.Item(16, 5).Value = 0.9890781214508737
$ws.Cells.Item(16, 6).Value = 1.035512109540261
$ws.Cells.Item(16, 9).Value = 1.032769449078452
$ws.Cells.Item(16, 10).Value = 1.033598783222627
$ws.Cells.Item(16, 11).Value = 1.034113205398479
$ws.Cells.Item(16, 12).Value = 0.9926916645766087
$ws.Cells.Item(16, 13).Value = 1.038951361844958
$ws.Cells.Item(16, 14).Value = 1.035066612050241

# Row 17
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.027888437389302
$ws.Cells.Item(17, 4).Value = 1.031153229434407
$ws.Cells.Item(17, 5).Value = 0.9894763578477731
$ws.Cells.Item(17, 6).Value = 1.036213137874287
$ws.Cells.Item(17, 9).Value = 1.032960179619202
$ws.Cells.Item(17, 10).Value = 1.03412501680262
$ws.Cells.Item(17, 11).Value = 1.034536048318171
$ws.Cells.Item(17, 12).Value = 0.9930127773692701
$ws.Cells.Item(17, 13).Value = 1.039578235570565
$ws.Cells.Item(17, 14).Value = 1.035593592942277

# Row 18
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.028277171787454
$ws.Cells.Item(18, 4).Value = 1.031442513770082
$ws.Cells.Item(18, 5).Value = 0.9897087662937551
$ws.Cells.Item(18, 6).Value = 1.036621597686448
$ws.Cells.Item(18, 9).Value = 1.033071032796045
$ws.Cells.Item(18, 10).Value = 1.034431457713723
$ws.Cells.Item(18, 11).Value = 1.034782197302424
$ws.Cells.Item(18, 12).Value = 0.9932001317071766
$ws.Cells.Item(18, 13).Value = 1.039943395815989
$ws.Cells.Item(18, 14).Value = 1.035900469034623

# Row 19
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.028409639004721
$ws.Cells.Item(19, 4).Value = 1.03154108036797
$ws.Cells.Item(19, 5).Value = 0.9897880325774039
$ws.Cells.Item(19, 6).Value = 1.036760798213312
$ws.Cells.Item(19, 9).Value = 1.033108763775688
$ws.Cells.Item(19, 10).Value = 1.034535861436034
$ws.Cells.Item(19, 11).Value = 1.034866045337919
$ws.Cells.Item(19, 12).Value = 0.993264023964098
$ws.Cells.Item(19, 13).Value = 1.040067824385646
$ws.Cells.Item(19, 14).Value = 1.036005021022195

# Row 20
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.027816893955387
$ws.Cells.Item(20, 4).Value = 1.031099983539949
$ws.Cells.Item(20, 5).Value = 0.9894336180355766
$ws.Cells.Item(20, 6).Value = 1.036137969634027
$ws.Cells.Item(20, 9).Value = 1.032939757120349
$ws.Cells.Item(20, 10).Value = 1.034068608982737
$ws.Cells.Item(20, 11).Value = 1.034490731858359
$ws.Cells.Item(20, 12).Value = 0.9929783193490043
$ws.Cells.Item(20, 13).Value = 1.039511028210014
$ws.Cells.Item(20, 14).Value = 1.03553710501682

# Row 21
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.025884112432458
$ws.Cells.Item(21, 4).Value = 1.029660894448305
$ws.Cells.Item(21, 5).Value = 0.9882828385668255
$ws.Cells.Item(21, 6).Value = 1.034107898625287
$ws.Cells.Item(21, 9).Value = 1.032385622255401
$ws.Cells.Item(21, 10).Value = 1.032543587374487
$ws.Cells.Item(21, 11).Value = 1.033264785065583
$ws.Cells.Item(21, 12).Value = 0.9920501090198107
$ws.Cells.Item(21, 13).Value = 1.037695092661263
$ws.Cells.Item(21, 14).Value = 1.034009917702963

# Row 22
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.024664416713426
$ws.Cells.Item(22, 4).Value = 1.028752147157914
$ws.Cells.Item(22, 5).Value = 0.9875604150241496
$ws.Cells.Item(22, 6).Value = 1.032827426118745
$ws.Cells.Item(22, 9).Value = 1.032033595237383
$ws.Cells.Item(22, 10).Value = 1.031580111311026
$ws.Cells.Item(22, 11).Value = 1.032489496574212
$ws.Cells.Item(22, 12).Value = 0.991467000034148
$ws.Cells.Item(22, 13).Value = 1.036548848078044
$ws.Cells.Item(22, 14).Value = 1.03304507339298

# Row 23
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.025311431745536
$ws.Cells.Item(23, 4).Value = 1.029234268904342
$ws.Cells.Item(23, 5).Value = 0.9879432794636459
$ws.Cells.Item(23, 6).Value = 1.033506622941453
$ws.Cells.Item(23, 9).Value = 1.032220556862449
$ws.Cells.Item(23, 10).Value = 1.032091312859189
$ws.Cells.Item(23, 11).Value = 1.032900921366602
$ws.Cells.Item(23, 12).Value = 0.9917760702887607
$ws.Cells.Item(23, 13).Value = 1.037156926240858
$ws.Cells.Item(23, 14).Value = 1.033557000905977

# Row 24
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.027849222842695
$ws.Cells.Item(24, 4).Value = 1.031124044384578
$ws.Cells.Item(24, 5).Value = 0.9894529299347241
$ws.Cells.Item(24, 6).Value = 1.036171936279713
$ws.Cells.Item(24, 9).Value = 1.032948986394242
$ws.Cells.Item(24, 10).Value = 1.034094098800077
$ws.Cells.Item(24, 11).Value = 1.034511209924021
$ws.Cells.Item(24, 12).Value = 0.9929938892766438
$ws.Cells.Item(24, 13).Value = 1.039541397820351
$ws.Cells.Item(24, 14).Value = 1.035562631032624

# Row 25
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.030775346904531
$ws.Cells.Item(25, 4).Value = 1.033300350254968
$ws.Cells.Item(25, 5).Value = 0.9912096547607046
$ws.Cells.Item(25, 6).Value = 1.03924776532729
$ws.Cells.Item(25, 9).Value = 1.033778748752085
$ws.Cells.Item(25, 10).Value = 1.036398573462286
$ws.Cells.Item(25, 11).Value = 1.036360756004403
$ws.Cells.Item(25, 12).Value = 0.9944092447426411
$ws.Cells.Item(25, 13).Value = 1.042289498540844
$ws.Cells.Item(25, 14).Value = 1.037870378313181
